# Revive the "CurrencyDataPoint" (AmountWithCurrency) component and use it
# only for the "Total Amount" (totalAmounts) rows of the EU Taxonomy
# Non-Financials framework data model sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

# Rows whose "Field Name" is "Absolute Share" (the totalAmounts fields):
# their Component column (F) moves from "Currency" to the new
# "AmountWithCurrency" component.
$totalAmountRows = @(12, 14, 16, 19, 31, 33, 35, 38, 50, 52, 54, 57)
foreach ($r in $totalAmountRows) {
    $ws.Cells.Item($r, 6).Value = "AmountWithCurrency"
}

# The two "Total Amount" group header rows get their Document-Support
# column (I) marked as "Extended".
$ws.Cells.Item(29, 9).Value = "Extended"
$ws.Cells.Item(48, 9).Value = "Extended"

# Minor view-state changes captured alongside the edit.
$excel.ActiveWindow.Zoom = 130
$null = $ws.Range("F10").Select()
